# Apply updated cryptocurrency price/volume data to Sheet1
# Commit: Updated symbol list on Mon Jan 23 10:18:51 UTC 2023 with GitHub Actions
#
# The Price (D) and Volume(1h) (E) columns in the source data are stored as
# plain text (e.g. "304.83", "0.99%") rather than real numbers, so we force
# the cells to Text number format before assigning the new values. This stops
# Excel's autodetection from silently turning "304.83" into the number
# 304.83 or "0.99%" into the fraction 0.0099 formatted as a percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D27","E27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "304.83"
$ws.Range("E2").Value = "0.99%"
$ws.Range("D3").Value = "35.64"
$ws.Range("E3").Value = "-4.55%"
$ws.Range("D4").Value = "5.089"
$ws.Range("E4").Value = "1.83%"
$ws.Range("D5").Value = "0.07870"
$ws.Range("E5").Value = "0.84%"
$ws.Range("D6").Value = "2.126"
$ws.Range("E6").Value = "-3.21%"
$ws.Range("D7").Value = "7.913"
$ws.Range("E7").Value = "-1.43%"
$ws.Range("D8").Value = "4.108"
$ws.Range("E8").Value = "1.78%"
$ws.Range("D9").Value = "0.9201"
$ws.Range("E9").Value = "0.61%"
$ws.Range("D10").Value = "0.09668"
$ws.Range("E10").Value = "0.03%"
$ws.Range("D11").Value = "0.1845"
$ws.Range("E11").Value = "-2.21%"
$ws.Range("D12").Value = "0.08627"
$ws.Range("E12").Value = "-1.34%"
$ws.Range("D13").Value = "0.03534"
$ws.Range("E13").Value = "0.16%"
$ws.Range("D14").Value = "0.09938"
$ws.Range("E14").Value = "-0.17%"
$ws.Range("D15").Value = "0.001447"
$ws.Range("E15").Value = "-2.05%"
$ws.Range("D16").Value = "0.005666"
$ws.Range("E16").Value = "-0.31%"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").Value = "-0.15%"
$ws.Range("D18").Value = "2.647"
$ws.Range("E18").Value = "11.31%"
$ws.Range("D19").Value = "0.3436"
$ws.Range("E19").Value = "-0.77%"
$ws.Range("E20").Value = "3.36%"
$ws.Range("D21").Value = "5.168"
$ws.Range("E21").Value = "8.28%"
$ws.Range("D22").Value = "0.2205"
$ws.Range("E22").Value = "-4.00%"
$ws.Range("D23").Value = "0.04543"
$ws.Range("E23").Value = "-1.95%"
$ws.Range("D24").Value = "0.005057"
$ws.Range("E24").Value = "5.63%"
$ws.Range("D25").Value = "0.001234"
$ws.Range("E25").Value = "0.20%"
$ws.Range("D27").Value = "0.0004751"
$ws.Range("E27").Value = "0.00%"
$ws.Range("D39").Value = "0.01838"
$ws.Range("E39").Value = "4.97%"
$ws.Range("D40").Value = "0.04735"
$ws.Range("E40").Value = "0.09%"
$ws.Range("D41").Value = "0.007513"
$ws.Range("E41").Value = "-6.48%"
$ws.Range("D42").Value = "0.1395"
$ws.Range("E42").Value = "0.41%"
$ws.Range("D43").Value = "0.007739"
$ws.Range("E43").Value = "0.63%"
$ws.Range("D44").Value = "0.002233"
$ws.Range("E44").Value = "0.43%"
$ws.Range("D45").Value = "0.01101"
$ws.Range("E45").Value = "5.96%"
$ws.Range("D46").Value = "0.00006322"
$ws.Range("E46").Value = "3.86%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.32%"
$ws.Range("D48").Value = "0.0005801"
$ws.Range("E48").Value = "0.01%"
$ws.Range("D49").Value = "47.51"
$ws.Range("E49").Value = "505.92%"
$ws.Range("D50").Value = "0.002001"
$ws.Range("E50").Value = "-25.64%"
$ws.Range("D51").Value = "0.00002096"
$ws.Range("E51").Value = "-0.32%"
